$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "season record" columns appended after the existing data (AC is the
# last currently-used column): Wins / Losses / Ties.

# Match the bold/centered/bordered header formatting used by the rest of
# row 1 (e.g. AC1) before writing the header text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Same record (96 wins, 66 losses, 0 ties) applied to every player row.
for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 30).Value = 96
    $ws.Cells.Item($r, 31).Value = 66
    $ws.Cells.Item($r, 32).Value = 0
}
